$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 31.14997866666667
$ws.Range("H2").Value = 93.44993600000001
$ws.Range("I2").Value = 0.4621739036316256
$ws.Range("J2").Value = 0.4621739036316256
$ws.Range("M2").Value = 20.88869433333333
$ws.Range("N2").Value = 62.666083
$ws.Range("O2").Value = 0.06073223131780172
$ws.Range("P2").Value = 0.06073223131780172
$ws.Range("Q2").Value = 650.6823828578542
$ws.Range("R2").Value = 5856.141445720688
$ws.Range("S2").Value = 0.02806885242440729
$ws.Range("T2").Value = 0.02806885242440729

$ws.Range("G3").Value = 31.14997866666667
$ws.Range("H3").Value = 93.44993600000001
$ws.Range("I3").Value = 0.4621739036316256
$ws.Range("J3").Value = 0.4621739036316256
$ws.Range("O3").Value = 0.1122209945765712
$ws.Range("P3").Value = 0.1122209945765712
$ws.Range("Q3").Value = 1202.330666490071
$ws.Range("R3").Value = 10820.97599841064
$ws.Range("S3").Value = 0.05186561513287741
$ws.Range("T3").Value = 0.05186561513287741

$ws.Range("G4").Value = 31.14997866666667
$ws.Range("H4").Value = 93.44993600000001
$ws.Range("I4").Value = 0.4621739036316256
$ws.Range("J4").Value = 0.4621739036316256
$ws.Range("M4").Value = 277.3327996666667
$ws.Range("N4").Value = 831.9983990000001
$ws.Range("O4").Value = 0.8063232422570387
$ws.Range("P4").Value = 0.8063232422570388
$ws.Range("Q4").Value = 8638.910793183608
$ws.Range("R4").Value = 77750.19713865248
$ws.Range("S4").Value = 0.3726615604628445
$ws.Range("T4").Value = 0.3726615604628445

$ws.Range("G5").Value = 31.14997866666667
$ws.Range("H5").Value = 93.44993600000001
$ws.Range("I5").Value = 0.4621739036316256
$ws.Range("J5").Value = 0.4621739036316256
$ws.Range("M5").Value = 7.127805333333332
$ws.Range("N5").Value = 21.383416
$ws.Range("O5").Value = 0.02072353184858837
$ws.Range("P5").Value = 0.02072353184858837
$ws.Range("Q5").Value = 222.0309840734862
$ws.Range("R5").Value = 1998.278856661376
$ws.Range("S5").Value = 0.009577875611496403
$ws.Range("T5").Value = 0.009577875611496405

$ws.Range("G6").Value = 18.94069966666667
$ws.Range("H6").Value = 56.822099
$ws.Range("I6").Value = 0.2810241764892454
$ws.Range("J6").Value = 0.2810241764892454
$ws.Range("M6").Value = 20.88869433333333
$ws.Range("N6").Value = 62.666083
$ws.Range("O6").Value = 0.06073223131780172
$ws.Range("P6").Value = 0.06073223131780172
$ws.Range("Q6").Value = 395.6464857964686
$ws.Range("R6").Value = 3560.818372168217
$ws.Range("S6").Value = 0.01706722529243958
$ws.Range("T6").Value = 0.01706722529243958

$ws.Range("G7").Value = 18.94069966666667
$ws.Range("H7").Value = 56.822099
$ws.Range("I7").Value = 0.2810241764892454
$ws.Range("J7").Value = 0.2810241764892454
$ws.Range("O7").Value = 0.1122209945765712
$ws.Range("P7").Value = 0.1122209945765712
$ws.Range("Q7").Value = 731.0754301857927
$ws.Range("R7").Value = 6579.678871672135
$ws.Range("S7").Value = 0.031536812585685
$ws.Range("T7").Value = 0.031536812585685

$ws.Range("G8").Value = 18.94069966666667
$ws.Range("H8").Value = 56.822099
$ws.Range("I8").Value = 0.2810241764892454
$ws.Range("J8").Value = 0.2810241764892454
$ws.Range("M8").Value = 277.3327996666667
$ws.Range("N8").Value = 831.9983990000001
$ws.Range("O8").Value = 0.8063232422570387
$ws.Range("P8").Value = 0.8063232422570388
$ws.Range("Q8").Value = 5252.877266202167
$ws.Range("R8").Value = 47275.89539581951
$ws.Range("S8").Value = 0.2265963251394226
$ws.Range("T8").Value = 0.2265963251394226

$ws.Range("G9").Value = 18.94069966666667
$ws.Range("H9").Value = 56.822099
$ws.Range("I9").Value = 0.2810241764892454
$ws.Range("J9").Value = 0.2810241764892454
$ws.Range("M9").Value = 7.127805333333332
$ws.Range("N9").Value = 21.383416
$ws.Range("O9").Value = 0.02072353184858837
$ws.Range("P9").Value = 0.02072353184858837
$ws.Range("Q9").Value = 135.0056201011315
$ws.Range("R9").Value = 1215.050580910184
$ws.Range("S9").Value = 0.005823813471698195
$ws.Range("T9").Value = 0.005823813471698196

$ws.Range("G10").Value = 14.86848
$ws.Range("H10").Value = 44.60544
$ws.Range("I10").Value = 0.2206044349565553
$ws.Range("J10").Value = 0.2206044349565553
$ws.Range("M10").Value = 20.88869433333333
$ws.Range("N10").Value = 62.666083
$ws.Range("O10").Value = 0.06073223131780172
$ws.Range("P10").Value = 0.06073223131780172
$ws.Range("Q10").Value = 310.58313392128
$ws.Range("R10").Value = 2795.24820529152
$ws.Range("S10").Value = 0.01339779957351446
$ws.Range("T10").Value = 0.01339779957351446

$ws.Range("G11").Value = 14.86848
$ws.Range("H11").Value = 44.60544
$ws.Range("I11").Value = 0.2206044349565553
$ws.Range("J11").Value = 0.2206044349565553
$ws.Range("O11").Value = 0.1122209945765712
$ws.Range("P11").Value = 0.1122209945765712
$ws.Range("Q11").Value = 573.8954000384
$ws.Range("R11").Value = 5165.0586003456
$ws.Range("S11").Value = 0.02475644909882715
$ws.Range("T11").Value = 0.02475644909882715

$ws.Range("G12").Value = 14.86848
$ws.Range("H12").Value = 44.60544
$ws.Range("I12").Value = 0.2206044349565553
$ws.Range("J12").Value = 0.2206044349565553
$ws.Range("M12").Value = 277.3327996666667
$ws.Range("N12").Value = 831.9983990000001
$ws.Range("O12").Value = 0.8063232422570387
$ws.Range("P12").Value = 0.8063232422570388
$ws.Range("Q12").Value = 4123.51718518784
$ws.Range("R12").Value = 37111.65466669056
$ws.Range("S12").Value = 0.1778784832504517
$ws.Range("T12").Value = 0.1778784832504517

$ws.Range("G13").Value = 14.86848
$ws.Range("H13").Value = 44.60544
$ws.Range("I13").Value = 0.2206044349565553
$ws.Range("J13").Value = 0.2206044349565553
$ws.Range("M13").Value = 7.127805333333332
$ws.Range("N13").Value = 21.383416
$ws.Range("O13").Value = 0.02072353184858837
$ws.Range("P13").Value = 0.02072353184858837
$ws.Range("Q13").Value = 105.97963104256
$ws.Range("R13").Value = 953.8166793830399
$ws.Range("S13").Value = 0.004571703033762015
$ws.Range("T13").Value = 0.004571703033762015

$ws.Range("G14").Value = 2.439668
$ws.Range("H14").Value = 7.319004
$ws.Range("I14").Value = 0.03619748492257375
$ws.Range("J14").Value = 0.03619748492257375
$ws.Range("M14").Value = 20.88869433333333
$ws.Range("N14").Value = 62.666083
$ws.Range("O14").Value = 0.06073223131780172
$ws.Range("P14").Value = 0.06073223131780172
$ws.Range("Q14").Value = 50.96147912681466
$ws.Range("R14").Value = 458.653312141332
$ws.Range("S14").Value = 0.002198354027440389
$ws.Range("T14").Value = 0.002198354027440389

$ws.Range("G15").Value = 2.439668
$ws.Range("H15").Value = 7.319004
$ws.Range("I15").Value = 0.03619748492257375
$ws.Range("J15").Value = 0.03619748492257375
$ws.Range("O15").Value = 0.1122209945765712
$ws.Range("P15").Value = 0.1122209945765712
$ws.Range("Q15").Value = 94.16660229027332
$ws.Range("R15").Value = 847.49942061246
$ws.Range("S15").Value = 0.004062117759181668
$ws.Range("T15").Value = 0.004062117759181668

$ws.Range("G16").Value = 2.439668
$ws.Range("H16").Value = 7.319004
$ws.Range("I16").Value = 0.03619748492257375
$ws.Range("J16").Value = 0.03619748492257375
$ws.Range("M16").Value = 277.3327996666667
$ws.Range("N16").Value = 831.9983990000001
$ws.Range("O16").Value = 0.8063232422570387
$ws.Range("P16").Value = 0.8063232422570388
$ws.Range("Q16").Value = 676.5999566971773
$ws.Range("R16").Value = 6089.399610274596
$ws.Range("S16").Value = 0.02918687340431994
$ws.Range("T16").Value = 0.02918687340431994

$ws.Range("G17").Value = 2.439668
$ws.Range("H17").Value = 7.319004
$ws.Range("I17").Value = 0.03619748492257375
$ws.Range("J17").Value = 0.03619748492257375
$ws.Range("M17").Value = 7.127805333333332
$ws.Range("N17").Value = 21.383416
$ws.Range("O17").Value = 0.02072353184858837
$ws.Range("P17").Value = 0.02072353184858837
$ws.Range("Q17").Value = 17.38947858196266
$ws.Range("R17").Value = 156.505307237664
$ws.Range("S17").Value = 0.0007501397316317542
$ws.Range("T17").Value = 0.0007501397316317543
